# Apply edits to benchmark-milestones worksheet:
#  - Rename several metric row labels in column A (rows 4-12)
#  - Flip a number of per-CVE result cells in rows 8-11 from "No" to "Yes"
#  - Append two new summary rows: "Test Iterations" (18) and "Number of Containers" (19)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename row labels in column A (rows 4-12) ---
$ws.Cells.Item(4, 1).Value = "cve_id_ok"
$ws.Cells.Item(5, 1).Value = "main_service"
$ws.Cells.Item(6, 1).Value = "main_version"
$ws.Cells.Item(7, 1).Value = "aux_services"
$ws.Cells.Item(8, 1).Value = "docker_runs"
$ws.Cells.Item(9, 1).Value = "services_ok"
$ws.Cells.Item(10, 1).Value = "code_main_version"
$ws.Cells.Item(11, 1).Value = "docker_vulnerable"
$ws.Cells.Item(12, 1).Value = "exploitable"

# --- Flip specific cells in rows 8-11 from "No" to "Yes" ---
$ws.Cells.Item(8, 8).Value = "Yes"  # H8
$ws.Cells.Item(8, 9).Value = "Yes"  # I8
$ws.Cells.Item(8, 10).Value = "Yes"  # J8
$ws.Cells.Item(8, 12).Value = "Yes"  # L8
$ws.Cells.Item(8, 14).Value = "Yes"  # N8
$ws.Cells.Item(8, 15).Value = "Yes"  # O8
$ws.Cells.Item(8, 16).Value = "Yes"  # P8
$ws.Cells.Item(8, 18).Value = "Yes"  # R8
$ws.Cells.Item(8, 20).Value = "Yes"  # T8
$ws.Cells.Item(8, 21).Value = "Yes"  # U8
$ws.Cells.Item(8, 22).Value = "Yes"  # V8
$ws.Cells.Item(8, 23).Value = "Yes"  # W8
$ws.Cells.Item(8, 24).Value = "Yes"  # X8
$ws.Cells.Item(8, 25).Value = "Yes"  # Y8
$ws.Cells.Item(8, 26).Value = "Yes"  # Z8
$ws.Cells.Item(8, 27).Value = "Yes"  # AA8
$ws.Cells.Item(8, 28).Value = "Yes"  # AB8
$ws.Cells.Item(8, 29).Value = "Yes"  # AC8
$ws.Cells.Item(8, 33).Value = "Yes"  # AG8
$ws.Cells.Item(8, 35).Value = "Yes"  # AI8
$ws.Cells.Item(8, 36).Value = "Yes"  # AJ8
$ws.Cells.Item(8, 37).Value = "Yes"  # AK8
$ws.Cells.Item(8, 38).Value = "Yes"  # AL8
$ws.Cells.Item(8, 39).Value = "Yes"  # AM8
$ws.Cells.Item(8, 42).Value = "Yes"  # AP8
$ws.Cells.Item(8, 43).Value = "Yes"  # AQ8
$ws.Cells.Item(8, 44).Value = "Yes"  # AR8
$ws.Cells.Item(8, 47).Value = "Yes"  # AU8
$ws.Cells.Item(8, 49).Value = "Yes"  # AW8
$ws.Cells.Item(8, 51).Value = "Yes"  # AY8
$ws.Cells.Item(8, 52).Value = "Yes"  # AZ8
$ws.Cells.Item(8, 53).Value = "Yes"  # BA8
$ws.Cells.Item(8, 54).Value = "Yes"  # BB8
$ws.Cells.Item(8, 55).Value = "Yes"  # BC8
$ws.Cells.Item(8, 56).Value = "Yes"  # BD8
$ws.Cells.Item(8, 57).Value = "Yes"  # BE8
$ws.Cells.Item(8, 58).Value = "Yes"  # BF8
$ws.Cells.Item(8, 59).Value = "Yes"  # BG8
$ws.Cells.Item(8, 60).Value = "Yes"  # BH8
$ws.Cells.Item(8, 61).Value = "Yes"  # BI8
$ws.Cells.Item(9, 10).Value = "Yes"  # J9
$ws.Cells.Item(9, 14).Value = "Yes"  # N9
$ws.Cells.Item(9, 15).Value = "Yes"  # O9
$ws.Cells.Item(9, 18).Value = "Yes"  # R9
$ws.Cells.Item(9, 21).Value = "Yes"  # U9
$ws.Cells.Item(9, 22).Value = "Yes"  # V9
$ws.Cells.Item(9, 23).Value = "Yes"  # W9
$ws.Cells.Item(9, 24).Value = "Yes"  # X9
$ws.Cells.Item(9, 25).Value = "Yes"  # Y9
$ws.Cells.Item(9, 26).Value = "Yes"  # Z9
$ws.Cells.Item(9, 27).Value = "Yes"  # AA9
$ws.Cells.Item(9, 28).Value = "Yes"  # AB9
$ws.Cells.Item(9, 29).Value = "Yes"  # AC9
$ws.Cells.Item(9, 33).Value = "Yes"  # AG9
$ws.Cells.Item(9, 36).Value = "Yes"  # AJ9
$ws.Cells.Item(9, 37).Value = "Yes"  # AK9
$ws.Cells.Item(9, 42).Value = "Yes"  # AP9
$ws.Cells.Item(9, 52).Value = "Yes"  # AZ9
$ws.Cells.Item(9, 53).Value = "Yes"  # BA9
$ws.Cells.Item(9, 54).Value = "Yes"  # BB9
$ws.Cells.Item(9, 55).Value = "Yes"  # BC9
$ws.Cells.Item(9, 56).Value = "Yes"  # BD9
$ws.Cells.Item(9, 57).Value = "Yes"  # BE9
$ws.Cells.Item(9, 58).Value = "Yes"  # BF9
$ws.Cells.Item(9, 59).Value = "Yes"  # BG9
$ws.Cells.Item(9, 61).Value = "Yes"  # BI9
$ws.Cells.Item(10, 8).Value = "Yes"  # H10
$ws.Cells.Item(10, 9).Value = "Yes"  # I10
$ws.Cells.Item(10, 10).Value = "Yes"  # J10
$ws.Cells.Item(10, 14).Value = "Yes"  # N10
$ws.Cells.Item(10, 15).Value = "Yes"  # O10
$ws.Cells.Item(10, 16).Value = "Yes"  # P10
$ws.Cells.Item(10, 18).Value = "Yes"  # R10
$ws.Cells.Item(10, 21).Value = "Yes"  # U10
$ws.Cells.Item(10, 22).Value = "Yes"  # V10
$ws.Cells.Item(10, 25).Value = "Yes"  # Y10
$ws.Cells.Item(10, 26).Value = "Yes"  # Z10
$ws.Cells.Item(10, 28).Value = "Yes"  # AB10
$ws.Cells.Item(10, 29).Value = "Yes"  # AC10
$ws.Cells.Item(10, 33).Value = "Yes"  # AG10
$ws.Cells.Item(10, 47).Value = "Yes"  # AU10
$ws.Cells.Item(10, 49).Value = "Yes"  # AW10
$ws.Cells.Item(10, 53).Value = "Yes"  # BA10
$ws.Cells.Item(10, 54).Value = "Yes"  # BB10
$ws.Cells.Item(10, 55).Value = "Yes"  # BC10
$ws.Cells.Item(10, 56).Value = "Yes"  # BD10
$ws.Cells.Item(10, 57).Value = "Yes"  # BE10
$ws.Cells.Item(10, 58).Value = "Yes"  # BF10
$ws.Cells.Item(10, 59).Value = "Yes"  # BG10
$ws.Cells.Item(10, 61).Value = "Yes"  # BI10
$ws.Cells.Item(11, 14).Value = "Yes"  # N11
$ws.Cells.Item(11, 15).Value = "Yes"  # O11
$ws.Cells.Item(11, 18).Value = "Yes"  # R11
$ws.Cells.Item(11, 26).Value = "Yes"  # Z11
$ws.Cells.Item(11, 27).Value = "Yes"  # AA11
$ws.Cells.Item(11, 38).Value = "Yes"  # AL11
$ws.Cells.Item(11, 42).Value = "Yes"  # AP11
$ws.Cells.Item(11, 44).Value = "Yes"  # AR11
$ws.Cells.Item(11, 59).Value = "Yes"  # BG11
$ws.Cells.Item(11, 61).Value = "Yes"  # BI11

# --- Add new row 18: "Test Iterations" ---
$ws.Cells.Item(18, 1).Value = "Test Iterations"
$ws.Cells.Item(18, 2).Value = 10
$ws.Cells.Item(18, 3).Value = 10
$ws.Cells.Item(18, 4).Value = 10
$ws.Cells.Item(18, 5).Value = 10
$ws.Cells.Item(18, 6).Value = 10
$ws.Cells.Item(18, 7).Value = 10
$ws.Cells.Item(18, 8).Value = 1
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 3
$ws.Cells.Item(18, 11).Value = 10
$ws.Cells.Item(18, 12).Value = 2
$ws.Cells.Item(18, 13).Value = 10
$ws.Cells.Item(18, 14).Value = 1
$ws.Cells.Item(18, 15).Value = 1
$ws.Cells.Item(18, 16).Value = 1
$ws.Cells.Item(18, 17).Value = 10
$ws.Cells.Item(18, 18).Value = 0
$ws.Cells.Item(18, 19).Value = 10
$ws.Cells.Item(18, 20).Value = 0
$ws.Cells.Item(18, 21).Value = 3
$ws.Cells.Item(18, 22).Value = 3
$ws.Cells.Item(18, 23).Value = 0
$ws.Cells.Item(18, 24).Value = 0
$ws.Cells.Item(18, 25).Value = 1
$ws.Cells.Item(18, 26).Value = 0
$ws.Cells.Item(18, 27).Value = 8
$ws.Cells.Item(18, 28).Value = 0
$ws.Cells.Item(18, 29).Value = 7
$ws.Cells.Item(18, 30).Value = 10
$ws.Cells.Item(18, 31).Value = 10
$ws.Cells.Item(18, 32).Value = 10
$ws.Cells.Item(18, 33).Value = 0
$ws.Cells.Item(18, 34).Value = 10
$ws.Cells.Item(18, 35).Value = 0
$ws.Cells.Item(18, 36).Value = 0
$ws.Cells.Item(18, 37).Value = 0
$ws.Cells.Item(18, 38).Value = 3
$ws.Cells.Item(18, 39).Value = 3
$ws.Cells.Item(18, 40).Value = 10
$ws.Cells.Item(18, 41).Value = 10
$ws.Cells.Item(18, 42).Value = 1
$ws.Cells.Item(18, 43).Value = 0
$ws.Cells.Item(18, 44).Value = 4
$ws.Cells.Item(18, 45).Value = 10
$ws.Cells.Item(18, 46).Value = 10
$ws.Cells.Item(18, 47).Value = 0
$ws.Cells.Item(18, 48).Value = 10
$ws.Cells.Item(18, 49).Value = 0
$ws.Cells.Item(18, 50).Value = 10
$ws.Cells.Item(18, 51).Value = 0
$ws.Cells.Item(18, 52).Value = 0
$ws.Cells.Item(18, 53).Value = 1
$ws.Cells.Item(18, 54).Value = 0
$ws.Cells.Item(18, 55).Value = 1
$ws.Cells.Item(18, 56).Value = 3
$ws.Cells.Item(18, 57).Value = 1
$ws.Cells.Item(18, 58).Value = 1
$ws.Cells.Item(18, 59).Value = 0
$ws.Cells.Item(18, 60).Value = 3
$ws.Cells.Item(18, 61).Value = 2

# --- Add new row 19: "Number of Containers" ---
$ws.Cells.Item(19, 1).Value = "Number of Containers"
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 2
$ws.Cells.Item(19, 9).Value = 2
$ws.Cells.Item(19, 10).Value = 2
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0
$ws.Cells.Item(19, 14).Value = 2
$ws.Cells.Item(19, 15).Value = 2
$ws.Cells.Item(19, 16).Value = 2
$ws.Cells.Item(19, 17).Value = 0
$ws.Cells.Item(19, 18).Value = 2
$ws.Cells.Item(19, 19).Value = 0
$ws.Cells.Item(19, 20).Value = 3
$ws.Cells.Item(19, 21).Value = 4
$ws.Cells.Item(19, 22).Value = 4
$ws.Cells.Item(19, 23).Value = 1
$ws.Cells.Item(19, 24).Value = 1
$ws.Cells.Item(19, 25).Value = 1
$ws.Cells.Item(19, 26).Value = 1
$ws.Cells.Item(19, 27).Value = 1
$ws.Cells.Item(19, 28).Value = 1
$ws.Cells.Item(19, 29).Value = 1
$ws.Cells.Item(19, 30).Value = 1
$ws.Cells.Item(19, 31).Value = 2
$ws.Cells.Item(19, 32).Value = 1
$ws.Cells.Item(19, 33).Value = 1
$ws.Cells.Item(19, 34).Value = 2
$ws.Cells.Item(19, 35).Value = 1
$ws.Cells.Item(19, 36).Value = 3
$ws.Cells.Item(19, 37).Value = 2
$ws.Cells.Item(19, 38).Value = 1
$ws.Cells.Item(19, 39).Value = 1
$ws.Cells.Item(19, 40).Value = 0
$ws.Cells.Item(19, 41).Value = 0
$ws.Cells.Item(19, 42).Value = 1
$ws.Cells.Item(19, 43).Value = 1
$ws.Cells.Item(19, 44).Value = 1
$ws.Cells.Item(19, 45).Value = 0
$ws.Cells.Item(19, 46).Value = 1
$ws.Cells.Item(19, 47).Value = 1
$ws.Cells.Item(19, 48).Value = 2
$ws.Cells.Item(19, 49).Value = 1
$ws.Cells.Item(19, 50).Value = 0
$ws.Cells.Item(19, 51).Value = 2
$ws.Cells.Item(19, 52).Value = 2
$ws.Cells.Item(19, 53).Value = 2
$ws.Cells.Item(19, 54).Value = 2
$ws.Cells.Item(19, 55).Value = 2
$ws.Cells.Item(19, 56).Value = 2
$ws.Cells.Item(19, 57).Value = 2
$ws.Cells.Item(19, 58).Value = 2
$ws.Cells.Item(19, 59).Value = 1
$ws.Cells.Item(19, 60).Value = 1
$ws.Cells.Item(19, 61).Value = 1

# --- Match the bold/bordered/centered style used by other row-label cells in column A ---
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

